$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 4235.5
$ws.Range("I82").Value = 471
$ws.Range("K82").Value = 1413
$ws.Range("M82").Value = -1007
$ws.Range("H85").Value = 4235.5
$ws.Range("I85").Value = 471
$ws.Range("K85").Value = 1413
$ws.Range("M85").Value = -9
$ws.Range("H98").Value = 2402.2778
$ws.Range("I98").Value = 2284.9092
$ws.Range("J98").Value = 3693.3333
$ws.Range("K98").Value = 2284.9092
$ws.Range("L98").Value = 3693.3333
$ws.Range("M98").Value = -786.9092000000001
$ws.Range("N98").Value = -6689.3333
$ws.Range("H112").Value = 1790
$ws.Range("J112").Value = 1288
$ws.Range("L112").Value = 3864
$ws.Range("N112").Value = -6080
$ws.Range("H113").Value = 4181.6313
$ws.Range("I113").Value = 3799.889
$ws.Range("K113").Value = 3799.889
$ws.Range("M113").Value = -545.8890000000001
$ws.Range("H122").Value = 2402.2778
$ws.Range("I122").Value = 2284.9092
$ws.Range("J122").Value = 3693.3333
$ws.Range("K122").Value = 6854.7276
$ws.Range("L122").Value = 11079.9999
$ws.Range("M122").Value = -4404.7276
$ws.Range("N122").Value = -15979.9999
$ws.Range("H138").Value = 4292.5127
$ws.Range("I138").Value = 2415.147
$ws.Range("J138").Value = 5680.1304
$ws.Range("K138").Value = 7245.441
$ws.Range("L138").Value = 17040.3912
$ws.Range("M138").Value = -2105.441
$ws.Range("N138").Value = -27320.3912

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 6961.6665
$ws.Range("I3").Value = 3192.5
$ws.Range("K3").Value = 3192.5
$ws.Range("M3").Value = -3077.5
$ws.Range("H22").Value = 1385.125
$ws.Range("I22").Value = 1385.125
$ws.Range("K22").Value = 1385.125
$ws.Range("M22").Value = -1086.125
$ws.Range("H32").Value = 18007.254
$ws.Range("I32").Value = 14997.182
$ws.Range("J32").Value = 26836.8
$ws.Range("K32").Value = 14997.182
$ws.Range("L32").Value = 26836.8
$ws.Range("M32").Value = -14710.182
$ws.Range("N32").Value = -27410.8
$ws.Range("H41").Value = 38907.89
$ws.Range("I41").Value = 1984.6
$ws.Range("J41").Value = 85062
$ws.Range("K41").Value = 1984.6
$ws.Range("L41").Value = 85062
$ws.Range("M41").Value = -1570.6
$ws.Range("N41").Value = -85890
$ws.Range("H45").Value = 1286.0233
$ws.Range("I45").Value = 1039.625
$ws.Range("J45").Value = 4571.3335
$ws.Range("K45").Value = 1039.625
$ws.Range("L45").Value = 4571.3335
$ws.Range("M45").Value = -662.625
$ws.Range("N45").Value = -5325.3335
$ws.Range("H61").Value = 2891.8823
$ws.Range("I61").Value = 2238.6667
$ws.Range("J61").Value = 3626.75
$ws.Range("K61").Value = 2238.6667
$ws.Range("L61").Value = 3626.75
$ws.Range("M61").Value = -2026.6667
$ws.Range("N61").Value = -4050.75
$ws.Range("H74").Value = 2131.0356
$ws.Range("I74").Value = 1648
$ws.Range("J74").Value = 3902.1667
$ws.Range("K74").Value = 1648
$ws.Range("L74").Value = 3902.1667
$ws.Range("M74").Value = -774
$ws.Range("N74").Value = -5650.1667
$ws.Range("H77").Value = 2131.0356
$ws.Range("I77").Value = 1648
$ws.Range("J77").Value = 3902.1667
$ws.Range("K77").Value = 8240
$ws.Range("L77").Value = 19510.8335
$ws.Range("M77").Value = -3872
$ws.Range("N77").Value = -28246.8335
$ws.Range("H92").Value = 34387.25
$ws.Range("J92").Value = 34387.25
$ws.Range("L92").Value = 34387.25
$ws.Range("N92").Value = -39379.25
$ws.Range("H122").Value = 2590.926
$ws.Range("I122").Value = 1785.6111
$ws.Range("K122").Value = 5356.8333
$ws.Range("M122").Value = -2906.8333
$ws.Range("H132").Value = 3658.6191
$ws.Range("I132").Value = 3073.4546
$ws.Range("J132").Value = 4302.3
$ws.Range("K132").Value = 9220.363799999999
$ws.Range("L132").Value = 12906.9
$ws.Range("M132").Value = -6690.363799999999
$ws.Range("N132").Value = -17966.9
$ws.Range("H133").Value = 27750
$ws.Range("J133").Value = 27750
$ws.Range("L133").Value = 27750
$ws.Range("N133").Value = -32810
$ws.Range("H136").Value = 2891.8823
$ws.Range("I136").Value = 2238.6667
$ws.Range("J136").Value = 3626.75
$ws.Range("K136").Value = 6716.000100000001
$ws.Range("L136").Value = 10880.25
$ws.Range("M136").Value = -4166.000100000001
$ws.Range("N136").Value = -15980.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 20910.416
$ws.Range("I82").Value = 4899.75
$ws.Range("J82").Value = 28915.75
$ws.Range("K82").Value = 4899.75
$ws.Range("L82").Value = 28915.75
$ws.Range("M82").Value = -4516.75
$ws.Range("N82").Value = -29681.75
$ws.Range("H85").Value = 20910.416
$ws.Range("I85").Value = 4899.75
$ws.Range("J85").Value = 28915.75
$ws.Range("K85").Value = 4899.75
$ws.Range("L85").Value = 28915.75
$ws.Range("M85").Value = -3573.75
$ws.Range("N85").Value = -31567.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1926401.5
$ws.Range("I31").Value = 2633557.5
$ws.Range("J31").Value = 6978.2856
$ws.Range("K31").Value = 2633557.5
$ws.Range("L31").Value = 6978.2856
$ws.Range("M31").Value = -2633262.5
$ws.Range("N31").Value = -7568.2856
$ws.Range("H34").Value = 1926401.5
$ws.Range("I34").Value = 2633557.5
$ws.Range("J34").Value = 6978.2856
$ws.Range("K34").Value = 2633557.5
$ws.Range("L34").Value = 6978.2856
$ws.Range("M34").Value = -2633355.5
$ws.Range("N34").Value = -7382.2856
$ws.Range("H62").Value = 3762.6667
$ws.Range("I62").Value = 2642.5
$ws.Range("J62").Value = 6003
$ws.Range("K62").Value = 2642.5
$ws.Range("L62").Value = 6003
$ws.Range("M62").Value = -2018.5
$ws.Range("N62").Value = -7251
$ws.Range("H65").Value = 3762.6667
$ws.Range("I65").Value = 2642.5
$ws.Range("J65").Value = 6003
$ws.Range("K65").Value = 13212.5
$ws.Range("L65").Value = 30015
$ws.Range("M65").Value = -10092.5
$ws.Range("N65").Value = -36255
$ws.Range("H96").Value = 23332.5
$ws.Range("J96").Value = 23332.5
$ws.Range("L96").Value = 23332.5
$ws.Range("N96").Value = -28824.5
$ws.Range("H123").Value = 30905.715
$ws.Range("J123").Value = 30905.715
$ws.Range("L123").Value = 30905.715
$ws.Range("N123").Value = -40705.715
$ws.Range("H141").Value = 30714.285
$ws.Range("I141").Value = 4750
$ws.Range("J141").Value = 32711.54
$ws.Range("K141").Value = 4750
$ws.Range("L141").Value = 32711.54
$ws.Range("M141").Value = 430
$ws.Range("N141").Value = -43071.54

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3125892.8
$ws.Range("I113").Value = 33333830
$ws.Range("J113").Value = 933.6896400000001
$ws.Range("K113").Value = 100001490
$ws.Range("L113").Value = 2801.06892
$ws.Range("M113").Value = -99999320
$ws.Range("N113").Value = -7141.06892
$ws.Range("H122").Value = 8768.615
$ws.Range("I122").Value = 584.1429000000001
$ws.Range("J122").Value = 18317.166
$ws.Range("K122").Value = 5257.2861
$ws.Range("L122").Value = 164854.494
$ws.Range("M122").Value = -2807.2861
$ws.Range("N122").Value = -169754.494
$ws.Range("H137").Value = 3266.95
$ws.Range("I137").Value = 2629.3333
$ws.Range("J137").Value = 5179.8
$ws.Range("K137").Value = 7887.999899999999
$ws.Range("L137").Value = 15539.4
$ws.Range("M137").Value = -2787.999899999999
$ws.Range("N137").Value = -25739.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 9474.75
$ws.Range("J92").Value = 9474.75
$ws.Range("L92").Value = 9474.75
$ws.Range("N92").Value = -13218.75
$ws.Range("I113").Value = 950
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 950
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = 1220
$ws.Range("N113").Value = -9840
$ws.Range("H132").Value = 4438.7104
$ws.Range("I132").Value = 5261.15
$ws.Range("J132").Value = 3524.889
$ws.Range("K132").Value = 15783.45
$ws.Range("L132").Value = 10574.667
$ws.Range("M132").Value = -13253.45
$ws.Range("N132").Value = -15634.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1056.8125
$ws.Range("I22").Value = 309.33334
$ws.Range("J22").Value = 1505.3
$ws.Range("K22").Value = 309.33334
$ws.Range("L22").Value = 1505.3
$ws.Range("M22").Value = -14.33334000000002
$ws.Range("N22").Value = -2095.3
$ws.Range("H27").Value = 1056.8125
$ws.Range("I27").Value = 309.33334
$ws.Range("J27").Value = 1505.3
$ws.Range("K27").Value = 309.33334
$ws.Range("L27").Value = 1505.3
$ws.Range("M27").Value = -202.33334
$ws.Range("N27").Value = -1719.3
$ws.Range("H46").Value = 1216.2069
$ws.Range("I46").Value = 895.38464
$ws.Range("J46").Value = 3996.6667
$ws.Range("K46").Value = 895.38464
$ws.Range("L46").Value = 3996.6667
$ws.Range("M46").Value = -707.38464
$ws.Range("N46").Value = -4372.6667
$ws.Range("H132").Value = 2527
$ws.Range("I132").Value = 1674.2858
$ws.Range("J132").Value = 4516.6665
$ws.Range("K132").Value = 5022.857400000001
$ws.Range("L132").Value = 13549.9995
$ws.Range("M132").Value = -2492.857400000001
$ws.Range("N132").Value = -18609.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 28700.334
$ws.Range("J103").Value = 28940.4
$ws.Range("L103").Value = 28940.4
$ws.Range("N103").Value = -31284.4
$ws.Range("H132").Value = 1733977
$ws.Range("I132").Value = 2002711.1
$ws.Range("K132").Value = 6008133.300000001
$ws.Range("M132").Value = -6005603.300000001
